$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Globo"
$ws.Range("B14").Value = "RJ TV 1"
$ws.Range("C14").Value = "Trânsito"
$ws.Range("D14").Value = "2025-04-01T12:52"
$ws.Range("E14").Value = "Neutro"
$ws.Range("F14").Value = "Caminhão invade o calçadão de Campos. Motorista teria errado o caminho e acabou subindo no calçadão do Centro. Repórter *ao vivo* do local. Imagens no Boulevard Francisco de Paula Carneiro, no Centro. Sem energia. Equipe no local"
